$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("PSFBeadsInput")
$wsInput.Range("A1").Value = "psf_bead_images"

$wsOutput = $wb.Worksheets.Item("PSFBeadsOutput")
$wsOutput.Range("I1").Value = "bead_properties"
$wsOutput.Range("J1").Value = "bead_z_profiles"
$wsOutput.Range("K1").Value = "bead_y_profiles"
$wsOutput.Range("L1").Value = "bead_x_profiles"
